$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4063.8408
$ws.Range("I15").Value = 4063.8408
$ws.Range("K15").Value = 12191.5224
$ws.Range("M15").Value = -12022.5224
$ws.Range("H74").Value = 6069
$ws.Range("I74").Value = 6033
$ws.Range("J74").Value = 6123
$ws.Range("K74").Value = 6033
$ws.Range("L74").Value = 6123
$ws.Range("M74").Value = -5097
$ws.Range("N74").Value = -7995
$ws.Range("H77").Value = 6069
$ws.Range("I77").Value = 6033
$ws.Range("J77").Value = 6123
$ws.Range("K77").Value = 30165
$ws.Range("L77").Value = 30615
$ws.Range("M77").Value = -25485
$ws.Range("N77").Value = -39975
$ws.Range("H100").Value = 1462.1904
$ws.Range("I100").Value = 1379
$ws.Range("J100").Value = 1728.4
$ws.Range("K100").Value = 1379
$ws.Range("L100").Value = 1728.4
$ws.Range("M100").Value = -838
$ws.Range("N100").Value = -2810.4
$ws.Range("H132").Value = 9319.054
$ws.Range("I132").Value = 6421.9614
$ws.Range("J132").Value = 16166.728
$ws.Range("K132").Value = 19265.8842
$ws.Range("L132").Value = 48500.18399999999
$ws.Range("M132").Value = -16735.8842
$ws.Range("N132").Value = -53560.18399999999
$ws.Range("H135").Value = 1298.1111
$ws.Range("I135").Value = 626.1429000000001
$ws.Range("K135").Value = 5635.2861
$ws.Range("M135").Value = -3100.2861
$ws.Range("H141").Value = 8673.4
$ws.Range("I141").Value = 9403.277
$ws.Range("K141").Value = 28209.831
$ws.Range("M141").Value = -23029.831
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 985.52
$ws.Range("I32").Value = 985.52
$ws.Range("K32").Value = 985.52
$ws.Range("M32").Value = -698.52
$ws.Range("H63").Value = 2403.6667
$ws.Range("I63").Value = 2244.5
$ws.Range("J63").Value = 3199.5
$ws.Range("K63").Value = 2244.5
$ws.Range("L63").Value = 3199.5
$ws.Range("M63").Value = -1558.5
$ws.Range("N63").Value = -4571.5
$ws.Range("H66").Value = 2403.6667
$ws.Range("I66").Value = 2244.5
$ws.Range("J66").Value = 3199.5
$ws.Range("K66").Value = 11222.5
$ws.Range("L66").Value = 15997.5
$ws.Range("M66").Value = -7790.5
$ws.Range("N66").Value = -22861.5
$ws.Range("H88").Value = 1512.6364
$ws.Range("I88").Value = 1462.8572
$ws.Range("K88").Value = 1462.8572
$ws.Range("M88").Value = -1056.8572
$ws.Range("H91").Value = 1512.6364
$ws.Range("I91").Value = 1462.8572
$ws.Range("K91").Value = 1462.8572
$ws.Range("M91").Value = -58.85719999999992
$ws.Range("H97").Value = 30334828
$ws.Range("J97").Value = 143870.14
$ws.Range("L97").Value = 143870.14
$ws.Range("N97").Value = -144862.14
$ws.Range("H102").Value = 3239.75
$ws.Range("I102").Value = 2916.8572
$ws.Range("J102").Value = 5500
$ws.Range("K102").Value = 2916.8572
$ws.Range("L102").Value = 5500
$ws.Range("M102").Value = -1294.8572
$ws.Range("N102").Value = -8744
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2564279.5
$ws.Range("I80").Value = 67.75
$ws.Range("J80").Value = 3703929.2
$ws.Range("K80").Value = 67.75
$ws.Range("L80").Value = 3703929.2
$ws.Range("M80").Value = 930.25
$ws.Range("N80").Value = -3705925.2
$ws.Range("H83").Value = 2564279.5
$ws.Range("I83").Value = 67.75
$ws.Range("J83").Value = 3703929.2
$ws.Range("K83").Value = 338.75
$ws.Range("L83").Value = 18519646
$ws.Range("M83").Value = 4653.25
$ws.Range("N83").Value = -18529630
$ws.Range("H86").Value = 2481.8235
$ws.Range("I86").Value = 1980.3
$ws.Range("J86").Value = 3198.2856
$ws.Range("K86").Value = 1980.3
$ws.Range("L86").Value = 3198.2856
$ws.Range("M86").Value = -857.3
$ws.Range("N86").Value = -5444.2856
$ws.Range("H89").Value = 2481.8235
$ws.Range("I89").Value = 1980.3
$ws.Range("J89").Value = 3198.2856
$ws.Range("K89").Value = 9901.5
$ws.Range("L89").Value = 15991.428
$ws.Range("M89").Value = -4285.5
$ws.Range("N89").Value = -27223.428
$ws.Range("H94").Value = 1127.4082
$ws.Range("I94").Value = 650.9697
$ws.Range("J94").Value = 2110.0625
$ws.Range("K94").Value = 650.9697
$ws.Range("L94").Value = 2110.0625
$ws.Range("M94").Value = -199.9697
$ws.Range("N94").Value = -3012.0625
$ws.Range("H99").Value = 1369.1923
$ws.Range("I99").Value = 1360.409
$ws.Range("J99").Value = 1417.5
$ws.Range("K99").Value = 1360.409
$ws.Range("L99").Value = 1417.5
$ws.Range("M99").Value = 137.5909999999999
$ws.Range("N99").Value = -4413.5
$ws.Range("H105").Value = 2667.1538
$ws.Range("I105").Value = 2652.0908
$ws.Range("K105").Value = 2652.0908
$ws.Range("M105").Value = -905.0907999999999
$ws.Range("H134").Value = 5224.25
$ws.Range("I134").Value = 4542.143
$ws.Range("K134").Value = 13626.429
$ws.Range("M134").Value = -11091.429
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 6879.9473
$ws.Range("I22").Value = 11345.182
$ws.Range("K22").Value = 11345.182
$ws.Range("M22").Value = -10995.182
$ws.Range("H29").Value = 9021
$ws.Range("J29").Value = 9021
$ws.Range("L29").Value = 9021
$ws.Range("N29").Value = -9607
$ws.Range("H35").Value = 1993
$ws.Range("I35").Value = 1993
$ws.Range("K35").Value = 1993
$ws.Range("M35").Value = -1699
$ws.Range("H36").Value = 15000
$ws.Range("I36").Value = 15000
$ws.Range("K36").Value = 15000
$ws.Range("M36").Value = -14612
$ws.Range("H40").Value = 15000
$ws.Range("I40").Value = 15000
$ws.Range("K40").Value = 15000
$ws.Range("M40").Value = -14840
$ws.Range("H58").Value = 1618.1951
$ws.Range("I58").Value = 1301.4412
$ws.Range("K58").Value = 1301.4412
$ws.Range("M58").Value = -1098.4412
$ws.Range("H136").Value = 1618.1951
$ws.Range("I136").Value = 1301.4412
$ws.Range("K136").Value = 3904.3236
$ws.Range("M136").Value = -1354.3236
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 2650.8572
$ws.Range("J52").Value = 2650.8572
$ws.Range("L52").Value = 7952.571599999999
$ws.Range("N52").Value = -8484.571599999999
$ws.Range("H80").Value = 3313.182
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3313.182
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 9939.545999999998
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -11811.546
$ws.Range("H83").Value = 3313.182
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3313.182
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 29818.638
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -39178.638
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 24004996
$ws.Range("I80").Value = 36927324
$ws.Range("K80").Value = 36927324
$ws.Range("M80").Value = -36926326
$ws.Range("H83").Value = 24004996
$ws.Range("I83").Value = 36927324
$ws.Range("K83").Value = 184636620
$ws.Range("M83").Value = -184631628
$ws.Range("H102").Value = 5091.9653
$ws.Range("J102").Value = 6755.1177
$ws.Range("L102").Value = 6755.1177
$ws.Range("N102").Value = -9999.117699999999
$ws.Range("H126").Value = 8667
$ws.Range("I126").Value = 8429.571
$ws.Range("J126").Value = 8999.4
$ws.Range("K126").Value = 25288.713
$ws.Range("L126").Value = 26998.2
$ws.Range("M126").Value = -22818.713
$ws.Range("N126").Value = -31938.2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 38332
$ws.Range("I42").Value = 39999
$ws.Range("J42").Value = 34998
$ws.Range("K42").Value = 39999
$ws.Range("L42").Value = 34998
$ws.Range("M42").Value = -39436
$ws.Range("N42").Value = -36124
$ws.Range("H49").Value = 38332
$ws.Range("I49").Value = 39999
$ws.Range("J49").Value = 34998
$ws.Range("K49").Value = 39999
$ws.Range("L49").Value = 34998
$ws.Range("M49").Value = -39852
$ws.Range("N49").Value = -35292
$ws.Range("H68").Value = 33335026
$ws.Range("I68").Value = 47620660
$ws.Range("K68").Value = 47620660
$ws.Range("M68").Value = -47619911
$ws.Range("H71").Value = 33335026
$ws.Range("I71").Value = 47620660
$ws.Range("K71").Value = 238103300
$ws.Range("M71").Value = -238099556
$ws.Range("H109").Value = 72299.5
$ws.Range("J109").Value = 72299.5
$ws.Range("L109").Value = 72299.5
$ws.Range("N109").Value = -75073.5
$ws.Range("H132").Value = 26499
$ws.Range("I132").Value = 13798.8
$ws.Range("K132").Value = 41396.39999999999
$ws.Range("M132").Value = -38866.39999999999
$ws.Range("H136").Value = 5195.7144
$ws.Range("I136").Value = 4669.222
$ws.Range("K136").Value = 14007.666
$ws.Range("M136").Value = -11457.666
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7580418.5
$ws.Range("I81").Value = 9094503
$ws.Range("K81").Value = 18189006
$ws.Range("M81").Value = -18187945
$ws.Range("H84").Value = 7580418.5
$ws.Range("I84").Value = 9094503
$ws.Range("K84").Value = 90945030
$ws.Range("M84").Value = -90939726
$ws.Range("H100").Value = 450.96155
$ws.Range("I100").Value = 479
$ws.Range("J100").Value = 296.75
$ws.Range("K100").Value = 958
$ws.Range("L100").Value = 593.5
$ws.Range("M100").Value = -417
$ws.Range("N100").Value = -1675.5
$ws.Range("H132").Value = 8825.272000000001
$ws.Range("I132").Value = 7574.9375
$ws.Range("J132").Value = 12159.5
$ws.Range("K132").Value = 22724.8125
$ws.Range("L132").Value = 36478.5
$ws.Range("M132").Value = -20194.8125
$ws.Range("N132").Value = -41538.5
